$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.316.48"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").Value = "3.099.38"
$ws.Range("E3").Value = "  +0.95%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.46"
$ws.Range("E5").Value = "  +0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.96"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "3.102.36"
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.19"
$ws.Range("E10").Value = "  -0.87%  "
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("E12").Value = "  +1.44%  "
$ws.Range("D13").Value = "3.640.11"
$ws.Range("E13").Value = "  +1.22%  "
$ws.Range("E14").Value = "  +0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.00"
$ws.Range("E15").Value = "  -0.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000163"
$ws.Range("E16").Value = "  -0.53%  "
$ws.Range("D17").Value = "57.447.29"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "3.101.94"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.09"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.04"
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "338.04"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.512"
$ws.Range("E24").Value = "  +2.21%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.65"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  +0.09%  "
$ws.Range("D28").Value = "0.0₃0910"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.49"
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.16"
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.86"
$ws.Range("E32").Value = "  +2.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.92"
$ws.Range("E33").Value = "  +0.94%  "
$ws.Range("E34").Value = "  -0.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "158.11"
$ws.Range("E35").Value = "  +1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.62"
$ws.Range("E36").Value = "  +1.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.10"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "27.12"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.27"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0657"
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("B41").Value = "RenzoRestakedETH"
$ws.Range("C41").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D41").Value = "3.146.55"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.687"
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.93"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.51"
$ws.Range("E44").Value = "  +10.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.79"
$ws.Range("E45").Value = "  +0.48%  "
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "2.304.11"
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  -0.11%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.971"
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.71"
$ws.Range("E50").Value = "  -1.37%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.00"
$ws.Range("E51").Value = "  +1.32%  "
